$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 781, pushing the existing rows 781-832
# down to 783-834 (formats/styles come along with the insert, matching
# how Excel copies the row above's formatting on EntireRow.Insert).
$ws.Rows.Item(781).EntireRow.Insert()
$ws.Rows.Item(781).EntireRow.Insert()

# New row 781 ("Primera" quality) for the newly reported week (2022-09-01 / serial 44826)
$ws.Range("A781").Value = 8
$ws.Range("B781").Value = "Terminal La Palmera de La Serena"
$ws.Range("C781").Value = "Coquimbo"
$ws.Range("D781").Value = 44826
$ws.Range("E781").Value = 4
$ws.Range("F781").Value = 100112023
$ws.Range("G781").Value = "Brócoli"
$ws.Range("H781").Value = "Sin especificar"
$ws.Range("I781").Value = "Primera"
$ws.Range("J781").Value = 2200
$ws.Range("K781").Value = 750
$ws.Range("L781").Value = 800
$ws.Range("M781").Value = 775
$ws.Range("N781").Value = "`$/unidad"
$ws.Range("O781").Value = "Provincia del Elquí"
$ws.Range("P781").Value = 775
$ws.Range("Q781").Value = 1
$ws.Range("R781").Value = "Hortaliza"

# New row 782 ("Segunda" quality) for the same new week
$ws.Range("A782").Value = 8
$ws.Range("B782").Value = "Terminal La Palmera de La Serena"
$ws.Range("C782").Value = "Coquimbo"
$ws.Range("D782").Value = 44826
$ws.Range("E782").Value = 4
$ws.Range("F782").Value = 100112023
$ws.Range("G782").Value = "Brócoli"
$ws.Range("H782").Value = "Sin especificar"
$ws.Range("I782").Value = "Segunda"
$ws.Range("J782").Value = 1340
$ws.Range("K782").Value = 650
$ws.Range("L782").Value = 700
$ws.Range("M782").Value = 675
$ws.Range("N782").Value = "`$/unidad"
$ws.Range("O782").Value = "Provincia del Elquí"
$ws.Range("P782").Value = 675
$ws.Range("Q782").Value = 1
$ws.Range("R782").Value = "Hortaliza"
